$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.548.51'
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").Value = '3.015.15'
$ws.Range("E3").Value = '  +2.23%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '''378.17'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("D6").Value = '''102.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = '''0.544'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.91%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +2.48%  '
$ws.Range("D10").Value = '''36.55'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.15%  '
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("E12").Value = '  +0.74%  '
$ws.Range("D13").Value = '3.494.41'
$ws.Range("E13").Value = '  +2.53%  '
$ws.Range("D14").Value = '''18.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.20%  '
$ws.Range("E15").Value = '  +1.21%  '
$ws.Range("D16").Value = '3.003.87'
$ws.Range("E16").Value = '  +2.14%  '
$ws.Range("E17").Value = '  -1.87%  '
$ws.Range("D18").Value = '''10.47'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -13.48%  '
$ws.Range("D19").Value = '51.528.52'
$ws.Range("E19").Value = '  +1.17%  '
$ws.Range("D20").Value = '''3.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("D23").Value = '''69.84'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.44%  '
$ws.Range("D24").Value = '''267.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("E25").Value = '  -3.66%  '
$ws.Range("D26").Value = '''8.21'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.67%  '
$ws.Range("D27").Value = '''7.51'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.55%  '
$ws.Range("E28").Value = '  +5.33%  '
$ws.Range("D30").Value = '''26.17'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.19%  '
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("D32").Value = '''10.29'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.53%  '
$ws.Range("D33").Value = '''34.04'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.62%  '
$ws.Range("D34").Value = '''50.72'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("D35").Value = '''0.0454'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.01%  '
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("E38").Value = '  +6.01%  '
$ws.Range("E39").Value = '  +3.66%  '
$ws.Range("D40").Value = '''0.284'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.29%  '
$ws.Range("D41").Value = '''2.60'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.43%  '
$ws.Range("E42").Value = '  +2.27%  '
$ws.Range("E43").Value = '  -0.21%  '
$ws.Range("D44").Value = '''3.73'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.55%  '
$ws.Range("D45").Value = '''122.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.12%  '
$ws.Range("D46").Value = '''21.50'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.82%  '
$ws.Range("D47").Value = '''2.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.16%  '
$ws.Range("E48").Value = '  +2.34%  '
$ws.Range("D49").Value = '2.030.92'
$ws.Range("E49").Value = '  +0.93%  '
$ws.Range("D50").Value = '3.313.24'
$ws.Range("E50").Value = '  +2.43%  '
$ws.Range("E51").Value = '  +1.31%  '
